$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "LR" path column (column B) from nicog's machine/paths to loren's.
$ws.Range("B2").Value = "C:\Users\loren\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase 3.8.2\MRSUT"
$ws.Range("B3").Value = "C:\Users\loren\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase 3.8.2\IOT"
$ws.Range("B4").Value = "C:\Users\loren\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase Hybrid 3.3.18"
$ws.Range("B5").Value = "C:\Users\loren\Documents\GitHub\MARIO Organization\GreenTechs\Database"
$ws.Range("B6").Value = "C:\Users\loren\Documents\GitHub\MARIO Organization\GreenTechs\Add sectors"
$ws.Range("B7").Value = "C:\Users\loren\Documents\GitHub\MARIO Organization\GreenTechs\Shocks"
$ws.Range("B8").Value = "C:\Users\loren\Documents\GitHub\MARIO Organization\GreenTechs\Results"
$ws.Range("B9").Value = "C:\Users\loren\Documents\GitHub\MARIO Organization\GreenTechs\Plots"
$ws.Range("B10").Value = "C:\Users\loren\Documents\GitHub\MARIO Organization\GreenTechs\Shocks\ShockMaster.xlsx"

# Move the active selection to B7, matching the saved cursor position.
$ws.Range("B7").Select()
